# "Ask a question changes"
# Updates the "Ask a question" row (row 18) on the Apis sheet:
#   - Api url (E18):      user/askQuestion          -> /questions/askAquestion
#   - params  (F18):      {userid:xxx,text:xxx,...} -> formatted JSON sample body
#                          (now wrapped, matching the multi-line text)
#   - response (H18):     {data:{updated},sucsess..} -> properly formatted JSON response
# Also restores the sheet's selection to G16 (it had been left on F18/A16-scrolled).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Apis")

# --- Api url ---
$ws.Range("E18").Value = "/questions/askAquestion"

# --- params (request body sample) ---
$paramsText = "{`n    ""userid"":6,""text"":""<h1>How do the college manage marks</h1> "",""tags"":[1,2,3]`n}"
$ws.Range("F18").Value = $paramsText
$ws.Range("F18").WrapText = $true

# --- response sample ---
$responseText = "{`n    ""data"": {`n        ""updated"": true`n    },`n    ""success"": true,`n    ""reason"": null`n}"
$ws.Range("H18").Value = $responseText

# --- restore view selection ---
$ws.Activate() | Out-Null
$ws.Range("G16").Select() | Out-Null
